$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (the source file keeps every D-column
# cell as a text cell, e.g. "246.88"), even though the text looks numeric.
# Force the Text number format on each D cell we touch before assigning the
# new value so Excel keeps storing it as text instead of auto-converting it
# to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "247.19"

$ws.Range("D4").Value = "5.453"

$ws.Range("D5").Value = "0.05686"

$ws.Range("D7").Value = "0.8002"

$ws.Range("D8").Value = "1.037"

$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1458"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.07295"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.03163"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.02938"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09277"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001654"
$ws.Range("E14").Value = "13BitForexTokenBF"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.215"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04708"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "0.0005891"
$ws.Range("E17").Value = "16OneONE"

$ws.Range("D18").Value = "0.006393"

$ws.Range("D19").Value = "0.005043"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"

$ws.Range("D20").Value = "0.001047"

$ws.Range("D23").Value = "3.772"

$ws.Range("D24").Value = "6.428"

$ws.Range("D25").Value = "2.126"

$ws.Range("D27").Value = "0.1279"

$ws.Range("D40").Value = "0.04089"

$ws.Range("D41").Value = "0.006919"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1043"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002970"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "0.008958"

$ws.Range("D45").Value = "0.00005833"

$ws.Range("D47").Value = "0.6826"

$ws.Range("D48").Value = "0.01003"
